$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 9377.929
$ws.Range("J17").Value = 9377.929
$ws.Range("L17").Value = 28133.787
$ws.Range("N17").Value = -28469.787

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 44931.5
$ws.Range("I98").Value = 1099
$ws.Range("K98").Value = 1099
$ws.Range("M98").Value = 399

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 44931.5
$ws.Range("I122").Value = 1099
$ws.Range("K122").Value = 3297
$ws.Range("M122").Value = -847

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2953.48
$ws.Range("I137").Value = 831.8570999999999
$ws.Range("J137").Value = 3778.5557
$ws.Range("K137").Value = 2495.5713
$ws.Range("L137").Value = 11335.6671
$ws.Range("M137").Value = 54.42870000000039
$ws.Range("N137").Value = -16435.6671

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1822.4086
$ws.Range("J138").Value = 2096.4583
$ws.Range("L138").Value = 6289.374899999999
$ws.Range("N138").Value = -16569.3749

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20309.967
$ws.Range("I32").Value = 19531.256
$ws.Range("J32").Value = 43671.332
$ws.Range("K32").Value = 19531.256
$ws.Range("L32").Value = 43671.332
$ws.Range("M32").Value = -19244.256
$ws.Range("N32").Value = -44245.332

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2840.625
$ws.Range("I45").Value = 2451.8333
$ws.Range("J45").Value = 4007
$ws.Range("K45").Value = 2451.8333
$ws.Range("L45").Value = 4007
$ws.Range("M45").Value = -2074.8333
$ws.Range("N45").Value = -4761

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2676.8206
$ws.Range("I61").Value = 1315.95
$ws.Range("K61").Value = 1315.95
$ws.Range("M61").Value = -1103.95

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3660
$ws.Range("I74").Value = 4173.2
$ws.Range("K74").Value = 4173.2
$ws.Range("M74").Value = -3299.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 3660
$ws.Range("I77").Value = 4173.2
$ws.Range("K77").Value = 20866
$ws.Range("M77").Value = -16498

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1718.6666
$ws.Range("I110").Value = 1207.8572
$ws.Range("K110").Value = 1207.8572
$ws.Range("M110").Value = 837.1428000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 9092854
$ws.Range("I132").Value = 16130366
$ws.Range("J132").Value = 2733.5833
$ws.Range("K132").Value = 48391098
$ws.Range("L132").Value = 8200.749899999999
$ws.Range("M132").Value = -48388568
$ws.Range("N132").Value = -13260.7499

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2676.8206
$ws.Range("I136").Value = 1315.95
$ws.Range("K136").Value = 3947.85
$ws.Range("M136").Value = -1397.85

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2154
$ws.Range("I99").Value = 2143.3333
$ws.Range("J99").Value = 2250
$ws.Range("K99").Value = 2143.3333
$ws.Range("L99").Value = 2250
$ws.Range("M99").Value = -645.3332999999998
$ws.Range("N99").Value = -5246

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1619.8055
$ws.Range("I107").Value = 1382.9
$ws.Range("K107").Value = 1382.9
$ws.Range("M107").Value = 537.0999999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3723.8572
$ws.Range("I134").Value = 2590.7083
$ws.Range("J134").Value = 4573.7188
$ws.Range("K134").Value = 7772.124899999999
$ws.Range("L134").Value = 13721.1564
$ws.Range("M134").Value = -5237.124899999999
$ws.Range("N134").Value = -18791.1564

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2488.02
$ws.Range("I31").Value = 1200.4791
$ws.Range("J31").Value = 3676.5193
$ws.Range("K31").Value = 1200.4791
$ws.Range("L31").Value = 3676.5193
$ws.Range("M31").Value = -905.4791
$ws.Range("N31").Value = -4266.5193

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2488.02
$ws.Range("I34").Value = 1200.4791
$ws.Range("J34").Value = 3676.5193
$ws.Range("K34").Value = 1200.4791
$ws.Range("L34").Value = 3676.5193
$ws.Range("M34").Value = -998.4791
$ws.Range("N34").Value = -4080.5193

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3481.6545
$ws.Range("J58").Value = 2649.2856
$ws.Range("L58").Value = 2649.2856
$ws.Range("N58").Value = -3055.2856

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2237.4666
$ws.Range("I99").Value = 1928
$ws.Range("J99").Value = 2350
$ws.Range("K99").Value = 1928
$ws.Range("L99").Value = 2350
$ws.Range("M99").Value = -430
$ws.Range("N99").Value = -5346

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 2237.4666
$ws.Range("I126").Value = 1928
$ws.Range("J126").Value = 2350
$ws.Range("K126").Value = 5784
$ws.Range("L126").Value = 7050
$ws.Range("M126").Value = -3314
$ws.Range("N126").Value = -11990

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 23327.893
$ws.Range("I132").Value = 1050.56
$ws.Range("K132").Value = 3151.68
$ws.Range("M132").Value = -621.6799999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1843.762
$ws.Range("I134").Value = 824.6
$ws.Range("J134").Value = 2770.2727
$ws.Range("K134").Value = 2473.8
$ws.Range("L134").Value = 8310.8181
$ws.Range("M134").Value = 61.19999999999982
$ws.Range("N134").Value = -13380.8181

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 3481.6545
$ws.Range("J136").Value = 2649.2856
$ws.Range("L136").Value = 7947.8568
$ws.Range("N136").Value = -13047.8568

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 6746
$ws.Range("I113").Value = 7813.75
$ws.Range("K113").Value = 7813.75
$ws.Range("M113").Value = -5643.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2124
$ws.Range("I122").Value = 1837.8125
$ws.Range("K122").Value = 5513.4375
$ws.Range("M122").Value = -3063.4375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4662.528
$ws.Range("I126").Value = 9540.308000000001
$ws.Range("K126").Value = 28620.924
$ws.Range("M126").Value = -26150.924

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 6411921
$ws.Range("J82").Value = 11906259
$ws.Range("L82").Value = 11906259
$ws.Range("N82").Value = -11906981

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 6411921
$ws.Range("J85").Value = 11906259
$ws.Range("L85").Value = 11906259
$ws.Range("N85").Value = -11908755

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2298.4324
$ws.Range("I136").Value = 1669.8214
$ws.Range("J136").Value = 4254.1113
$ws.Range("K136").Value = 5009.4642
$ws.Range("L136").Value = 12762.3339
$ws.Range("M136").Value = -2459.4642
$ws.Range("N136").Value = -17862.3339

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 534.17645
$ws.Range("I100").Value = 433
$ws.Range("K100").Value = 866
$ws.Range("M100").Value = -325

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2598491
$ws.Range("I122").Value = 4082486.8
$ws.Range("J122").Value = 1498.75
$ws.Range("K122").Value = 12247460.4
$ws.Range("L122").Value = 4496.25
$ws.Range("M122").Value = -12245010.4
$ws.Range("N122").Value = -9396.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1841321.1
$ws.Range("I126").Value = 1964042.5
$ws.Range("K126").Value = 5892127.5
$ws.Range("M126").Value = -5889657.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1515.8334
$ws.Range("I132").Value = 1072.0667
$ws.Range("J132").Value = 2255.4443
$ws.Range("K132").Value = 3216.2001
$ws.Range("L132").Value = 6766.3329
$ws.Range("M132").Value = -686.2001
$ws.Range("N132").Value = -11826.3329

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 15537.684
$ws.Range("I136").Value = 29013.514
$ws.Range("K136").Value = 87040.542
$ws.Range("M136").Value = -84490.542
